$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 2030
    $ws.Range("F4").Value = 244
    $ws.Range("F6").Value = 6361
    $ws.Range("F7").Value = 235
    $ws.Range("F8").Value = 117
}
